# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect newly scraped numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1233
$ws1.Range("F5").Value = 2140
$ws1.Range("F8").Value = 133
$ws1.Range("F10").Value = 498
$ws1.Range("F16").Value = 7417
$ws1.Range("F21").Value = 1725
$ws1.Range("F33").Value = 70
$ws1.Range("F35").Value = 3905

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1233
$ws4.Range("F10").Value = 2140
$ws4.Range("F14").Value = 133
$ws4.Range("F18").Value = 498
$ws4.Range("F24").Value = 7417
$ws4.Range("F29").Value = 1725
$ws4.Range("F43").Value = 3905
